$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 95-96), pushing the
# existing rows 95-217 down to 97-219. This mirrors the weekly refresh where
# two new price records are prepended to the dataset.
$ws.Range("A95:A96").EntireRow.Insert()

# New row 95: Crespo record / Primera, 2022-04-20
$ws.Range("A95").Value = 7
$ws.Range("B95").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C95").Value = "Ñuble"
$ws.Range("D95").Value = 44671
$ws.Range("E95").Value = 16
$ws.Range("F95").Value = 100112006
$ws.Range("G95").Value = "Repollo"
$ws.Range("H95").Value = "Crespo record"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 120
$ws.Range("K95").Value = 1000
$ws.Range("L95").Value = 1100
$ws.Range("M95").Value = 1050
$ws.Range("N95").Value = "$/unidad"
$ws.Range("O95").Value = "Provincia de Diguillín"
$ws.Range("P95").Value = 1050
$ws.Range("Q95").Value = 1
$ws.Range("R95").Value = "Hortaliza"

# New row 96: Crespo record / Segunda, 2022-04-20
$ws.Range("A96").Value = 7
$ws.Range("B96").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C96").Value = "Ñuble"
$ws.Range("D96").Value = 44671
$ws.Range("E96").Value = 16
$ws.Range("F96").Value = 100112006
$ws.Range("G96").Value = "Repollo"
$ws.Range("H96").Value = "Crespo record"
$ws.Range("I96").Value = "Segunda"
$ws.Range("J96").Value = 80
$ws.Range("K96").Value = 900
$ws.Range("L96").Value = 900
$ws.Range("M96").Value = 900
$ws.Range("N96").Value = "$/unidad"
$ws.Range("O96").Value = "Provincia de Diguillín"
$ws.Range("P96").Value = 900
$ws.Range("Q96").Value = 1
$ws.Range("R96").Value = "Hortaliza"
